$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Microcontrolador (ESP32) / qty 1 ---
$ws.Range("A2").Style = "Normal"
$ws.Range("A2").Value = "Microcontrolador (ESP32)"
$ws.Range("A2").VerticalAlignment = -4108
$ws.Range("A2").HorizontalAlignment = -4131
$ws.Range("A2").WrapText = $true
$ws.Range("B2").Value = 1

# --- Row 3: Módulo GPS (u-blox NEO-6M) / qty 2 (no wrap, just centered) ---
$ws.Range("A3").Style = "Normal"
$ws.Range("A3").Value = "Módulo GPS (u-blox NEO-6M)"
$ws.Range("A3").VerticalAlignment = -4108
$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("B3").Value = 2

# --- Row 4: Módulo GSM (SIM800L) / qty 3 -- reuse A2's format ---
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = "Módulo GSM (SIM800L)"
$ws.Range("B4").Value = 3

# --- Row 5: Antena GPS / qty 4 (new row) -- reuse A2's format ---
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "Antena GPS"
$ws.Range("B5").Value = 4

$ws.Application.CutCopyMode = $false

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 24.65
$ws.Columns.Item(2).ColumnWidth = 10.0

# --- Data validation on B2:B1048576 (input message only, no restriction) ---
$r = $ws.Range("B2:B1048576")
$r.Validation.Add(0, 1, 1, "")
$r.Validation.ErrorTitle = "Erro de validação!"
$r.Validation.ErrorMessage = "O dado inserido não corresponde a um valor decimal."
$r.Validation.InputMessage = "Use a virgula para valores decimais."
$r.Validation.ShowInput = $true
$r.Validation.ShowError = $true

# --- Selection moves to B6 ---
$ws.Range("B6").Select()
